$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-27 Friday", "2024-09-28 Saturday"),
    @("72×92=6624", "40×67=2680"),
    @("69×33=2277", "77×72=5544"),
    @("34×34=1156", "86×83=7138"),
    @("72×70=5040", "64×11=704"),
    @("99×39=3861", "76×42=3192"),
    @("75×34=2550", "22×26=572"),
    @("32×20=640", "94×70=6580"),
    @("65×70=4550", "70×92=6440"),
    @("86×69=5934", "57×43=2451"),
    @("51×29=1479", "56×97=5432"),
    @("65×49=3185", "80×61=4880"),
    @("52×69=3588", "35×21=735"),
    @("87×35=3045", "19×59=1121"),
    @("90×49=4410", "63×78=4914"),
    @("63×42=2646", "47×53=2491"),
    @("74×46=3404", "39×57=2223"),
    @("36×34=1224", "53×21=1113"),
    @("27×35=945", "56×91=5096"),
    @("79×52=4108", "43×14=602"),
    @("52×70=3640", "53×20=1060"),
    @("56×82=4592", "12×14=168"),
    @("64×15=960", "39×23=897"),
    @("47×43=2021", "50×22=1100"),
    @("11×19=209", "27×45=1215"),
    @("92×23=2116", "69×84=5796")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
